# Add a new worksheet "Primos Serene 002" to the workbook, populate it with
# a second cash-invoice table, and make it the active sheet (mirroring the
# first sheet's layout/styles).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add the new worksheet right after the existing one ----------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Primos Serene 002"

# --- Column widths (match the look of sheet1, but per target dims) -----
$ws2.Columns.Item(1).ColumnWidth = 6.6640625
$ws2.Columns.Item(2).ColumnWidth = 33

# --- Header row ----------------------------------------------------------
$ws2.Range("A1").Value = "Sr. No"
$ws2.Range("B1").Value = "Item Description"
$ws2.Range("C1").Value = "Qty"
$ws2.Range("D1").Value = "Price"
$ws2.Range("E1").Value = "Total"
$ws2.Range("A1:E1").Style = $ws1.Range("A1:E1").Style
$ws2.Rows.Item(1).RowHeight = 18.6

# --- Data rows -------------------------------------------------------
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = "2 MP CP Plus Bullet Camera"
$ws2.Range("C2").Value = 1
$ws2.Range("D2").Value = 2700
$ws2.Range("E2").Formula = "=C2*D2"

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = "Camera Box"
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = 250
$ws2.Range("E3").Formula = "=C3*D3"

$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = "BNC Connector"
$ws2.Range("C4").Value = 2
$ws2.Range("D4").Value = 75
$ws2.Range("E4").Formula = "=C4*D4"

$ws2.Range("A5").Value = 4
$ws2.Range("B5").Value = "Service Call Charges"
$ws2.Range("C5").Value = 1
$ws2.Range("D5").Value = 750
$ws2.Range("E5").Formula = "=C5*D5"
# (kept in this order intentionally: camera, box, connector, service call)

$ws2.Range("A2:E5").Style = $ws1.Range("A2:E2").Style
$ws2.Rows.Item(2).RowHeight = 12.6

# --- Totals row ----------------------------------------------------------
$ws2.Range("A6").Value = "Total"
$ws2.Range("E6").Formula = "=SUM(E2:E5)"
$ws2.Range("A6:D6").Style = $ws1.Range("A3:D3").Style
$ws2.Range("E6").Style = $ws1.Range("E3").Style
$ws2.Range("A6:D6").Merge()

# --- Selection / view state ---------------------------------------------
$ws2.Range("Q5").Select()

# The new sheet becomes the active tab; the original sheet should no
# longer be flagged as the selected tab.
$ws2.Activate()
